$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 318793.47
$ws.Range("I15").Value = 318793.47
$ws.Range("K15").Value = 956380.4099999999
$ws.Range("M15").Value = -956211.4099999999

$ws.Range("H62").Value = 14364.9
$ws.Range("I62").Value = 22080.4
$ws.Range("J62").Value = 6649.4
$ws.Range("K62").Value = 22080.4
$ws.Range("L62").Value = 6649.4
$ws.Range("M62").Value = -21456.4
$ws.Range("N62").Value = -7897.4

$ws.Range("H65").Value = 14364.9
$ws.Range("I65").Value = 22080.4
$ws.Range("J65").Value = 6649.4
$ws.Range("K65").Value = 110402
$ws.Range("L65").Value = 33247
$ws.Range("M65").Value = -107282
$ws.Range("N65").Value = -39487

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H136").Value = 28947
$ws.Range("J136").Value = 28947
$ws.Range("L136").Value = 28947
$ws.Range("N136").Value = -39147

$ws.Range("H139").Value = 34710
$ws.Range("J139").Value = 34710
$ws.Range("L139").Value = 34710
$ws.Range("N139").Value = -44990

$ws.Range("H141").Value = 2756.4285
$ws.Range("I141").Value = 859
$ws.Range("J141").Value = 7500
$ws.Range("K141").Value = 2577
$ws.Range("L141").Value = 22500
$ws.Range("M141").Value = 2603
$ws.Range("N141").Value = -32860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2098.3
$ws.Range("I2").Value = 2595
$ws.Range("J2").Value = 1353.25
$ws.Range("K2").Value = 2595
$ws.Range("L2").Value = 1353.25
$ws.Range("M2").Value = -2482
$ws.Range("N2").Value = -1579.25

$ws.Range("H32").Value = 7913.5415
$ws.Range("I32").Value = 5134.6924
$ws.Range("J32").Value = 19955.223
$ws.Range("K32").Value = 5134.6924
$ws.Range("L32").Value = 19955.223
$ws.Range("M32").Value = -4847.6924
$ws.Range("N32").Value = -20529.223

$ws.Range("H45").Value = 1548.3636
$ws.Range("I45").Value = 1491.5
$ws.Range("J45").Value = 1700
$ws.Range("K45").Value = 1491.5
$ws.Range("L45").Value = 1700
$ws.Range("M45").Value = -1114.5
$ws.Range("N45").Value = -2454

$ws.Range("H74").Value = 264075.28
$ws.Range("I74").Value = 346287.88
$ws.Range("J74").Value = 80678
$ws.Range("K74").Value = 346287.88
$ws.Range("L74").Value = 80678
$ws.Range("M74").Value = -345413.88
$ws.Range("N74").Value = -82426

$ws.Range("H77").Value = 264075.28
$ws.Range("I77").Value = 346287.88
$ws.Range("J77").Value = 80678
$ws.Range("K77").Value = 1731439.4
$ws.Range("L77").Value = 403390
$ws.Range("M77").Value = -1727071.4
$ws.Range("N77").Value = -412126

$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988

$ws.Range("H116").Value = 2098.3
$ws.Range("I116").Value = 2595
$ws.Range("J116").Value = 1353.25
$ws.Range("K116").Value = 2595
$ws.Range("L116").Value = 1353.25
$ws.Range("M116").Value = -301
$ws.Range("N116").Value = -5941.25

$ws.Range("H132").Value = 3068.1404
$ws.Range("I132").Value = 3248.6206
$ws.Range("J132").Value = 2881.2144
$ws.Range("K132").Value = 9745.861800000001
$ws.Range("L132").Value = 8643.643199999999
$ws.Range("M132").Value = -7215.861800000001
$ws.Range("N132").Value = -13703.6432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2098.3
$ws.Range("I3").Value = 2595
$ws.Range("J3").Value = 1353.25
$ws.Range("K3").Value = 2595
$ws.Range("L3").Value = 1353.25
$ws.Range("M3").Value = -2481
$ws.Range("N3").Value = -1581.25

$ws.Range("H64").Value = 421.42105
$ws.Range("J64").Value = 435.2
$ws.Range("L64").Value = 435.2
$ws.Range("N64").Value = -885.2

$ws.Range("H67").Value = 421.42105
$ws.Range("J67").Value = 435.2
$ws.Range("L67").Value = 435.2
$ws.Range("N67").Value = -1995.2

$ws.Range("H99").Value = 4044.7144
$ws.Range("I99").Value = 4435.933
$ws.Range("J99").Value = 3066.6667
$ws.Range("K99").Value = 4435.933
$ws.Range("L99").Value = 3066.6667
$ws.Range("M99").Value = -2937.933
$ws.Range("N99").Value = -6062.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2830.3225
$ws.Range("I31").Value = 1566.9231
$ws.Range("J31").Value = 9400
$ws.Range("K31").Value = 1566.9231
$ws.Range("L31").Value = 9400
$ws.Range("M31").Value = -1271.9231
$ws.Range("N31").Value = -9990

$ws.Range("H34").Value = 2830.3225
$ws.Range("I34").Value = 1566.9231
$ws.Range("J34").Value = 9400
$ws.Range("K34").Value = 1566.9231
$ws.Range("L34").Value = 9400
$ws.Range("M34").Value = -1364.9231
$ws.Range("N34").Value = -9804

$ws.Range("H94").Value = 3955.077
$ws.Range("I94").Value = 706.2857
$ws.Range("J94").Value = 7745.3335
$ws.Range("K94").Value = 706.2857
$ws.Range("L94").Value = 7745.3335
$ws.Range("M94").Value = -255.2857
$ws.Range("N94").Value = -8647.333500000001

$ws.Range("H99").Value = 64091.188
$ws.Range("I99").Value = 92210.45
$ws.Range("J99").Value = 2228.8
$ws.Range("K99").Value = 92210.45
$ws.Range("L99").Value = 2228.8
$ws.Range("M99").Value = -90712.45
$ws.Range("N99").Value = -5224.8

$ws.Range("H126").Value = 64091.188
$ws.Range("I126").Value = 92210.45
$ws.Range("J126").Value = 2228.8
$ws.Range("K126").Value = 276631.35
$ws.Range("L126").Value = 6686.400000000001
$ws.Range("M126").Value = -274161.35
$ws.Range("N126").Value = -11626.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 437.42856
$ws.Range("J97").Value = 460.33334
$ws.Range("L97").Value = 1381.00002
$ws.Range("N97").Value = -2373.00002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3943.3962
$ws.Range("I80").Value = 4460.6484
$ws.Range("J80").Value = 2747.25
$ws.Range("K80").Value = 4460.6484
$ws.Range("L80").Value = 2747.25
$ws.Range("M80").Value = -3462.6484
$ws.Range("N80").Value = -4743.25

$ws.Range("H83").Value = 3943.3962
$ws.Range("I83").Value = 4460.6484
$ws.Range("J83").Value = 2747.25
$ws.Range("K83").Value = 22303.242
$ws.Range("L83").Value = 13736.25
$ws.Range("M83").Value = -17311.242
$ws.Range("N83").Value = -23720.25
